$wb = $excel.ActiveWorkbook

# --- Challenge 3 worksheet: insert two new "Alerting and monitoring" rows ---
$ws3 = $wb.Worksheets.Item("Challenge 3")

# Insert two new blank rows above the existing row 20 ("Log Analytics and
# Security Center..."), shifting everything below down by two rows and
# expanding the section SUM() formulas automatically.
$ws3.Rows.Item(20).EntireRow.Insert()
$ws3.Rows.Item(20).EntireRow.Insert()

# Match the row height used by the neighboring detail rows.
$ws3.Rows.Item(20).RowHeight = 17
$ws3.Rows.Item(21).RowHeight = 17

# Populate the two new rows with their objective text.
$ws3.Cells.Item(20, 2).Value = "Alerts have been configured to report approaching system thresholds"
$ws3.Cells.Item(21, 2).Value = "Alerts have been configured to report underutilized resources"

# Update the "Alerting and monitoring" section's Points Possible totals.
$ws3.Cells.Item(16, 4).Value = 36
$ws3.Cells.Item(16, 5).Value = 29

# Make Challenge 3 the active sheet/tab, with E17 as the selected cell.
$ws3.Activate()
$ws3.Range("E17").Select()
